# Trade #55 closed at 2026-02-18 00:21:09 - unknown UNKNOWN +0.000%
#
# This records the closing of the open "momentum" trade (row 85 on "All
# Trades" / row 15 on "momentum") via an early_exit, and the subsequent
# opening of a brand-new "momentum" trade (appended as trade #113).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - headline counters
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 83      # Total Trades
$summary.Range("B9").Value = 50.6    # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - "momentum" strategy row (row 11)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D11").Value = 13      # Trades
$status.Range("G11").Value = 15.38   # Win Rate %

# ---------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# 3a. Close out trade #84 (sheet row 85)
$allTrades.Range("G85").Value = 0.87
$allTrades.Range("H85").Value = "CLOSED"
$allTrades.Range("K85").Value = 99.70999999999999
$allTrades.Range("L85").Value = "early_exit"
$allTrades.Range("M85").Value = 0.12

# 3b. Append the newly opened trade #113 (sheet row 114)
$allTrades.Range("A114").Value = 113
$allTrades.Range("B114").Value = "'2026-02-18"
$allTrades.Range("C114").Value = "00:21:03"
$allTrades.Range("D114").Value = "momentum"
$allTrades.Range("E114").Value = "DOWN"
$allTrades.Range("F114").Value = 0.87
$allTrades.Range("H114").Value = "OPEN"
$allTrades.Range("I114").Value = 0
$allTrades.Range("J114").Value = 0
$allTrades.Range("K114").Value = 99.7087371310913
$allTrades.Range("M114").Value = 0
$allTrades.Range("N114").Value = 0
$allTrades.Range("O114").Value = 0
$allTrades.Range("P114").Value = 0.9
$allTrades.Range("Q114").Value = "Downward momentum: -3.810% over 10 samples"

# ---------------------------------------------------------------------
# 4. "momentum" strategy sheet
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

# 4a. Close out trade #84 (sheet row 15)
$momentum.Range("G15").Value = 0.87
$momentum.Range("H15").Value = "CLOSED"
$momentum.Range("K15").Value = 99.70999999999999
$momentum.Range("P15").Value = "early_exit"
$momentum.Range("Q15").Value = 0.12

# 4b. Append the newly opened trade #113 (sheet row 32)
$momentum.Range("A32").Value = 113
$momentum.Range("B32").Value = "'2026-02-18"
$momentum.Range("C32").Value = "00:21:03"
$momentum.Range("D32").Value = "momentum"
$momentum.Range("E32").Value = "DOWN"
$momentum.Range("F32").Value = 0.87
$momentum.Range("H32").Value = "OPEN"
$momentum.Range("I32").Value = 0
$momentum.Range("J32").Value = 0
$momentum.Range("K32").Value = 99.7087371310913
$momentum.Range("L32").Value = 0
$momentum.Range("M32").Value = 0
$momentum.Range("N32").Value = 0.9
$momentum.Range("O32").Value = "Downward momentum: -3.810% over 10 samples"
$momentum.Range("Q32").Value = 0
